$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Silvia%Maffoni%NULL%1,         Silvia%Brazzo%NULL%1,         Rachele%De Giuseppe%NULL%1,         Ginevra%Biino%NULL%2,         Ginevra%Biino%NULL%0,         Ilaria%Vietti%NULL%1,         Cristina%Pallavicini%NULL%1,         Hellas%Cena%NULL%1,         Emma%Ruiz Moreno%NULL%3,         Emma%Ruiz Moreno%NULL%0,         Emma%Ruiz Moreno%NULL%0]'
$ws.Range("E3").Value = '[Jaideep%khare%NULL%1,         Sushil%jindal%NULL%1]'
$ws.Range("E4").Value = '[Sithum%Munasinghe%NULL%1,         Sandro%Sperandei%NULL%1,         Louise%Freebairn%NULL%1,         Elizabeth%Conroy%NULL%1,         Hir%Jani%NULL%1,         Sandra%Marjanovic%NULL%1,         Andrew%Page%NULL%1]'
$ws.Range("E5").Value = '[Rachel G.%Curtis%NULL%1,         Timothy%Olds%NULL%2,         Timothy%Olds%NULL%0,         Ty%Ferguson%NULL%1,         François%Fraysse%NULL%2,         François%Fraysse%NULL%0,         Dorothea%Dumuid%NULL%1,         Adrian%Esterman%NULL%1,         Gilly A.%Hendrie%NULL%1,         Wendy J.%Brown%NULL%1,         Rajini%Lagiseti%NULL%1,         Carol A.%Maher%NULL%1,         Stephanie%Prince Ware%NULL%2,         Stephanie%Prince Ware%NULL%0]'
$ws.Range("E6").Value = '[Idoia%Imaz-Aramburu%NULL%1,         Ana-Belén%Fraile-Bermúdez%NULL%2,         Ana-Belén%Fraile-Bermúdez%NULL%0,         Batirtze San%Martín-Gamboa%NULL%2,         Batirtze San%Martín-Gamboa%NULL%0,         Silvia%Cepeda-Miguel%NULL%1,         Borja%Doncel-García%NULL%2,         Borja%Doncel-García%NULL%0,         Ainhoa%Fernandez-Atutxa%NULL%2,         Ainhoa%Fernandez-Atutxa%NULL%0,         Amaia%Irazusta%NULL%1,         Idoia%Zarrazquin%NULL%2,         Idoia%Zarrazquin%NULL%0,         Simona%Bo%NULL%4,         Simona%Bo%NULL%0,         Simona%Bo%NULL%0]'
$ws.Range("E7").Value = '[Arenaza%Lide%coreGivesNoEmail%1,        Cadenas-S\u00e1nchez%Cristina%coreGivesNoEmail%1,        Labayen%Go\u00f1i Idoia%coreGivesNoEmail%1,        Medrano%Echeverr\u00eda Mar\u00eda%coreGivesNoEmail%1,        Os\u00e9s%Recalde Maddi%coreGivesNoEmail%1]'
$ws.Range("E8").Value = '[Felix%Naughton%f.naughton@uea.ac.uk%1,         Emma%Ward%NULL%2,         Emma%Ward%NULL%0,         Mizanur%Khondoker%NULL%1,         Pippa%Belderson%NULL%1,         Anne%Marie Minihane%NULL%1,         Jack%Dainty%NULL%1,         Sarah%Hanson%NULL%1,         Richard%Holland%NULL%1,         Tracey%Brown%NULL%1,         Caitlin%Notley%NULL%1]'
$ws.Range("E9").Value = '[Moritz%Herle%NULL%1,         Andrea D.%Smith%NULL%1,         Feifei%Bu%NULL%1,         Andrew%Steptoe%NULL%1,         Daisy%Fancourt%d.fancourt@ucl.ac.uk%1]'
$ws.Range("E10").Value = '[Yukako%Hosomi%NULL%1,         Chihiro%Munekawa%NULL%1,         Yoshitaka%Hashimoto%y-hashi@koto.kpu-m.ac.jp%1,         Takuro%Okamura%NULL%2,         Takuro%Okamura%NULL%0,         Fuyuko%Takahashi%NULL%1,         Rena%Kawano%NULL%1,         Hanako%Nakajima%NULL%1,         Saori%Majima%NULL%1,         Takafumi%Senmaru%NULL%1,         Naoko%Nakanishi%NULL%1,         Emi%Ushigome%NULL%1,         Masahide%Hamaguchi%NULL%1,         Masahiro%Yamazaki%NULL%1,         Michiaki%Fukui%NULL%1]'
$ws.Range("E11").Value = '[Koryu%Sato%NULL%1,         Satomi%Kobayashi%NULL%1,         Mai%Yamaguchi%NULL%1,         Ryohei%Sakata%NULL%1,         Yuki%Sasaki%NULL%1,         Chiaki%Murayama%NULL%1,         Naoki%Kondo%NULL%1]'
$ws.Range("E12").Value = '[Angelo%Pietrobelli%NULL%1,         Luca%Pecoraro%NULL%1,         Alessandro%Ferruzzi%NULL%1,         Moonseong%Heo%NULL%1,         Myles%Faith%NULL%1,         Thomas%Zoller%NULL%2,         Thomas%Zoller%NULL%0,         Franco%Antoniazzi%NULL%1,         Giorgio%Piacentini%NULL%1,         S. Nicole%Fearnbach%NULL%1,         Steven B.%Heymsfield%steven.heymsfield@pbrc.edu%2,         Steven B.%Heymsfield%steven.heymsfield@pbrc.edu%0]'
$ws.Range("E13").Value = '[Caruso%I.%coreGivesNoEmail%1,        Cignarelli%A.%coreGivesNoEmail%1,        Di%Molfetta S.%coreGivesNoEmail%1,        Giordano%F.%coreGivesNoEmail%1,        Giorgino%F.%coreGivesNoEmail%1,        Guarini%F.%coreGivesNoEmail%1,        Laviola%L.%coreGivesNoEmail%1,        Leonardini%A.%coreGivesNoEmail%1,        Natalicchio%A.%coreGivesNoEmail%1,        Perrini%S.%coreGivesNoEmail%1]'
$ws.Range("E14").Value = '[Marianna%Pellegrini%NULL%1,         Valentina%Ponzo%NULL%1,         Rosalba%Rosato%NULL%1,         Elena%Scumaci%NULL%1,         Ilaria%Goitre%NULL%1,         Andrea%Benso%NULL%1,         Sara%Belcastro%NULL%1,         Chiara%Crespi%NULL%1,         Franco%De Michieli%NULL%1,         Ezio%Ghigo%NULL%1,         Fabio%Broglio%NULL%1,         Simona%Bo%NULL%0]'
